$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 923 (weekly update: new week's
# prices inserted, all subsequent rows shift down by 3).
$ws.Rows("923:925").Insert()

# New row 923 - Mango, Especial, Peru, 2022-02-18
$ws.Cells.Item(923,1).Value = 6
$ws.Cells.Item(923,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(923,3).Value = "Metropolitana"
$ws.Cells.Item(923,4).Value = "2022-02-18"
$ws.Cells.Item(923,5).Value = 13
$ws.Cells.Item(923,6).Value = "Fruta"
$ws.Cells.Item(923,7).Value = 100108
$ws.Cells.Item(923,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(923,9).Value = 100108002
$ws.Cells.Item(923,10).Value = "Mango"
$ws.Cells.Item(923,11).Value = "Sin especificar"
$ws.Cells.Item(923,12).Value = "Especial"
$ws.Cells.Item(923,13).Value = 456
$ws.Cells.Item(923,14).Value = 5000
$ws.Cells.Item(923,15).Value = 6000
$ws.Cells.Item(923,16).Value = 5500
$ws.Cells.Item(923,17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(923,18).Value = "Perú"
$ws.Cells.Item(923,19).Value = 1375
$ws.Cells.Item(923,20).Value = 4

# New row 924 - Mango, Primera, Peru, 2022-02-18
$ws.Cells.Item(924,1).Value = 6
$ws.Cells.Item(924,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(924,3).Value = "Metropolitana"
$ws.Cells.Item(924,4).Value = "2022-02-18"
$ws.Cells.Item(924,5).Value = 13
$ws.Cells.Item(924,6).Value = "Fruta"
$ws.Cells.Item(924,7).Value = 100108
$ws.Cells.Item(924,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(924,9).Value = 100108002
$ws.Cells.Item(924,10).Value = "Mango"
$ws.Cells.Item(924,11).Value = "Sin especificar"
$ws.Cells.Item(924,12).Value = "Primera"
$ws.Cells.Item(924,13).Value = 456
$ws.Cells.Item(924,14).Value = 5000
$ws.Cells.Item(924,15).Value = 6000
$ws.Cells.Item(924,16).Value = 5500
$ws.Cells.Item(924,17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(924,18).Value = "Perú"
$ws.Cells.Item(924,19).Value = 1375
$ws.Cells.Item(924,20).Value = 4

# New row 925 - Mango, Segunda, Peru, 2022-02-18
$ws.Cells.Item(925,1).Value = 6
$ws.Cells.Item(925,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(925,3).Value = "Metropolitana"
$ws.Cells.Item(925,4).Value = "2022-02-18"
$ws.Cells.Item(925,5).Value = 13
$ws.Cells.Item(925,6).Value = "Fruta"
$ws.Cells.Item(925,7).Value = 100108
$ws.Cells.Item(925,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(925,9).Value = 100108002
$ws.Cells.Item(925,10).Value = "Mango"
$ws.Cells.Item(925,11).Value = "Sin especificar"
$ws.Cells.Item(925,12).Value = "Segunda"
$ws.Cells.Item(925,13).Value = 456
$ws.Cells.Item(925,14).Value = 5000
$ws.Cells.Item(925,15).Value = 6000
$ws.Cells.Item(925,16).Value = 5500
$ws.Cells.Item(925,17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(925,18).Value = "Perú"
$ws.Cells.Item(925,19).Value = 1375
$ws.Cells.Item(925,20).Value = 4
